# 1 April - Friday, heading home.
# The "departure date" originally entered in C2 (Mexico City trip) is moved
# down to row 10, and C2 is updated to a new origin city (London) to match
# the new search (no longer a date, so it loses the date-ish style).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("flightsSearch")

# Preserve C2's current value + style by copying it down to the new C10
# cell before we overwrite C2.
$ws.Range("C2").Copy($ws.Range("C10"))

# C2 becomes "London" (same text already used in B2) with default
# (unstyled) formatting.
$ws.Range("C2").ClearFormats()
$ws.Range("C2").Value = "London"

# Move the active selection to C2, matching the new point of interest.
[void]$ws.Range("C2").Select()
